# 485-RBI-EPP-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-EarlyRePayment-Loanproduct4.xlsx
# "code refactoring and loan accounting and charges added"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoan_Input
$ws2 = $wb.Worksheets.Item(2)   # ProductLoan_Output

# ---------------------------------------------------------------------------
# 1. Field value corrections on ProductLoan_Input
# ---------------------------------------------------------------------------

# shortname: was the text "kar2" -> now the number 485
$ws1.Range("B3").Value = 485

# nominalinterestratedefault: 12 -> 1
$ws1.Range("B11").Value = 1

# maximumallowedaoutstandingbalance: 5000 -> 10000
$ws1.Range("B26").Value = 10000

# ---------------------------------------------------------------------------
# 2. New loan-accounting / charges rows appended below the existing data
#    (rows 29-40), re-using the same row style as the other label/value
#    pairs (A column style copied from A13, B column style copied from B13)
# ---------------------------------------------------------------------------

$newRows = @(
    @("fundsource", "Cash"),
    @("loanprotfolio", "Loan portfolio "),
    @("interestreceivable", "Interest Receivable "),
    @("penaltiesreceivable", "Penalties Receivable "),
    @("transferinsuspense", "Transfer in Suspence "),
    @("feesreceivable", "Fees Receivable"),
    @("incomefrominterest", "Income from interest"),
    @("incomefrompenalties", "Income from penalties"),
    @("incomefromfees", "Income from fees"),
    @("incomefromrecoveryrepayments", "Income from recovery repayments"),
    @("loseswrittenoff", "Losses Writtenoff "),
    @("overpaymentliability", "Overpayment Liability")
)

$r = 29
foreach ($pair in $newRows) {
    $ws1.Range("A13").Copy()
    $ws1.Range("A$r").PasteSpecial(-4122)
    $ws1.Range("B13").Copy()
    $ws1.Range("B$r").PasteSpecial(-4122)

    $ws1.Range("A$r").Value = $pair[0]
    $ws1.Range("B$r").Value = $pair[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Column B on ProductLoan_Input gets noticeably wider (and loses the
#    "best fit" auto flag in favour of an explicit custom width)
# ---------------------------------------------------------------------------

$ws1.Columns.Item(2).ColumnWidth = 56.6

# ---------------------------------------------------------------------------
# 4. Selection / scroll bookkeeping to match the new, longer sheet
# ---------------------------------------------------------------------------

$ws2.Range("B1").Select()
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws1.Range("A41").Select()

# ---------------------------------------------------------------------------
# 5. ProductLoan_Output!B1 picks up a fresh (general-alignment) style
# ---------------------------------------------------------------------------

$ws2.Range("B1").HorizontalAlignment = 1

Write-Output "edit complete"
